# Weekly update: insert two new price records (rows 318-319) for
# "Vega Modelo de Temuco - Perejil", shifting the existing historic
# rows down by two (318->320 ... 399->401).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 318; formatting (incl. the date number format
# on column D) is inherited from the row above, matching the rest of
# the table.
$ws.Rows("318:319").Insert()

# New row 318
$ws.Range("A318").Value = 10
$ws.Range("B318").Value = "Vega Modelo de Temuco"
$ws.Range("C318").Value = "La Araucanía"
$ws.Range("D318").Value = 44855
$ws.Range("E318").Value = 9
$ws.Range("F318").Value = 100112044
$ws.Range("G318").Value = "Perejil"
$ws.Range("H318").Value = "Sin especificar"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 50
$ws.Range("K318").Value = 4000
$ws.Range("L318").Value = 4000
$ws.Range("M318").Value = 4000
$ws.Range("N318").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O318").Value = "Provincia de Cautín"
$ws.Range("P318").Value = 1333
$ws.Range("Q318").Value = 3
$ws.Range("R318").Value = "Hortaliza"

# New row 319
$ws.Range("A319").Value = 10
$ws.Range("B319").Value = "Vega Modelo de Temuco"
$ws.Range("C319").Value = "La Araucanía"
$ws.Range("D319").Value = 44855
$ws.Range("E319").Value = 9
$ws.Range("F319").Value = 100112044
$ws.Range("G319").Value = "Perejil"
$ws.Range("H319").Value = "Sin especificar"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 50
$ws.Range("K319").Value = 3000
$ws.Range("L319").Value = 3000
$ws.Range("M319").Value = 3000
$ws.Range("N319").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O319").Value = "Región Metropolitana"
$ws.Range("P319").Value = 1000
$ws.Range("Q319").Value = 3
$ws.Range("R319").Value = "Hortaliza"
